$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37
$ws.Range("D37").Value = 44755
$ws.Range("K37").Value = 16000
$ws.Range("L37").Value = 17000
$ws.Range("M37").Value = 16500
$ws.Range("P37").Value = 330

# Row 38
$ws.Range("D38").Value = 44755
$ws.Range("K38").Value = 18000
$ws.Range("M38").Value = 19000
$ws.Range("P38").Value = 633

# Row 39
$ws.Range("D39").Value = 44364
$ws.Range("H39").Value = 'Argentina(o)'
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 19000
$ws.Range("L39").Value = 20000
$ws.Range("M39").Value = 19500
$ws.Range("N39").Value = '$/caja 50 unidades'
$ws.Range("P39").Value = 390
$ws.Range("Q39").Value = 50

# Row 40
$ws.Range("D40").Value = 44364
$ws.Range("H40").Value = 'Española'
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 19000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = 19500
$ws.Range("N40").Value = '$/caja 30 unidades'
$ws.Range("P40").Value = 650
$ws.Range("Q40").Value = 30

# Row 41
$ws.Range("D41").Value = 44727
$ws.Range("J41").Value = 150
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 21000
$ws.Range("M41").Value = 20467
$ws.Range("P41").Value = 682

# Row 42
$ws.Range("D42").Value = 44748
$ws.Range("J42").Value = 110
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 16000
$ws.Range("M42").Value = 15545
$ws.Range("N42").Value = '$/caja 40 unidades'
$ws.Range("P42").Value = 389
$ws.Range("Q42").Value = 40

# Row 43
$ws.Range("D43").Value = 44699
$ws.Range("K43").Value = 19000
$ws.Range("M43").Value = 19500
$ws.Range("P43").Value = 650

# Row 44
$ws.Range("D44").Value = 44358
$ws.Range("H44").Value = 'Argentina(o)'
$ws.Range("K44").Value = 18000
$ws.Range("L44").Value = 20000
$ws.Range("M44").Value = 19000
$ws.Range("N44").Value = '$/caja 50 unidades'
$ws.Range("P44").Value = 380
$ws.Range("Q44").Value = 50

# Row 45
$ws.Range("D45").Value = 44358
$ws.Range("H45").Value = 'Española'
$ws.Range("K45").Value = 18000
$ws.Range("L45").Value = 20000
$ws.Range("M45").Value = 19000
$ws.Range("N45").Value = '$/caja 30 unidades'
$ws.Range("P45").Value = 633
$ws.Range("Q45").Value = 30

# Row 46
$ws.Range("D46").Value = 44428
$ws.Range("K46").Value = 14000
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = 14500
$ws.Range("P46").Value = 483

# Row 47
$ws.Range("D47").Value = 44489
$ws.Range("H47").Value = 'Argentina(o)'
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 10000
$ws.Range("M47").Value = 9500
$ws.Range("N47").Value = '$/caja 50 unidades'
$ws.Range("P47").Value = 190
$ws.Range("Q47").Value = 50

# Row 48
$ws.Range("D48").Value = 44489
$ws.Range("K48").Value = 8000
$ws.Range("L48").Value = 8500
$ws.Range("M48").Value = 8250
$ws.Range("P48").Value = 275

# Row 49
$ws.Range("D49").Value = 44426
$ws.Range("J49").Value = 50
$ws.Range("K49").Value = 12000
$ws.Range("L49").Value = 13000
$ws.Range("M49").Value = 12600
$ws.Range("O49").Value = 'Provincia de Limarí'
$ws.Range("P49").Value = 315

# Row 50 (new)
$ws.Range("A50").Value = 11
$ws.Range("B50").Value = 'Vega Monumental Concepción'
$ws.Range("C50").Value = 'Bíobío'
$ws.Range("D50").Value = 44376
$ws.Range("E50").Value = 8
$ws.Range("F50").Value = 100112013
$ws.Range("G50").Value = 'Alcachofa'
$ws.Range("H50").Value = 'Española'
$ws.Range("I50").Value = 'Primera'
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 19000
$ws.Range("L50").Value = 20000
$ws.Range("M50").Value = 19500
$ws.Range("N50").Value = '$/caja 30 unidades'
$ws.Range("O50").Value = 'Provincia de Limarí'
$ws.Range("P50").Value = 650
$ws.Range("Q50").Value = 30
$ws.Range("R50").Value = 'Hortaliza'
$ws.Range("D50").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 51 (new)
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = 'Vega Monumental Concepción'
$ws.Range("C51").Value = 'Bíobío'
$ws.Range("D51").Value = 44454
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112013
$ws.Range("G51").Value = 'Alcachofa'
$ws.Range("H51").Value = 'Madrigal'
$ws.Range("I51").Value = 'Primera'
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 13000
$ws.Range("L51").Value = 14000
$ws.Range("M51").Value = 13500
$ws.Range("N51").Value = '$/caja 40 unidades'
$ws.Range("O51").Value = 'Provincia del Elquí'
$ws.Range("P51").Value = 338
$ws.Range("Q51").Value = 40
$ws.Range("R51").Value = 'Hortaliza'
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Done"